# Applies the "Fixed workload plan and added the final plan" edit:
#  - Fills in the team members' student ids (row 5) and names (row 6)
#  - Fills in the actual workload percentages for each task (was all zeros)
#  - Updates the active selection on the "workload" sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("workload")

# --- Team member identification (rows 5 & 6) ---
# member 1 (D) / member 2 (E) / member 3 (F)
$ws.Range("D5").Value = 5762340
$ws.Range("E5").Value = 5556910
$ws.Range("F5").Value = 5699193

# Note: write F6/E6/D6 in this order so new shared strings are appended
# in the same order as the reference workbook (Andrea, Eduard, Tudor).
$ws.Range("F6").Value = "Andrea Vezzuto"
$ws.Range("E6").Value = "Eduard Faraon"
$ws.Range("D6").Value = "Tudor Coman"

# --- Basic features workload split (rows 8-15) ---
$ws.Range("D8").Value = 100
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0

$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 100
$ws.Range("F9").Value = 0

$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 100
$ws.Range("F10").Value = 0

$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 100

$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 100
$ws.Range("F12").Value = 0

$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 100
$ws.Range("F13").Value = 0

$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 100

$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 100
$ws.Range("F15").Value = 0

# --- Extra features workload split (rows 19-24) ---
$ws.Range("D19").Value = 100
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 0

$ws.Range("D20").Value = 100
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0

$ws.Range("D21").Value = 25
$ws.Range("E21").Value = 25
$ws.Range("F21").Value = 50

$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 100

$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 100

$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 100
$ws.Range("F24").Value = 0

# --- Update the active selection / active cell shown on the sheet ---
$ws.Range("L35").Select() | Out-Null
